$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
Write-Host $ws.Name
